$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Delete row 9 (K21-64-FI#6); subsequent rows (old row 10) shift up
$ws.Rows.Item(9).Delete()

# Step 2: Update filename (column B) values with the _CRR_DiadFit suffix where applicable
$ws.Range("B2").Value = "K21-59-FI#1_CRR_DiadFit"
$ws.Range("B3").Value = "K21-59-FI#2"
$ws.Range("B4").Value = "K21-59-FI#3_CRR_DiadFit"
$ws.Range("B5").Value = "K21-61-FI#1-rep1_CRR_DiadFit"
$ws.Range("B6").Value = "K21-61-FI#1-rep2_CRR_DiadFit"
$ws.Range("B7").Value = "K21-63-FI#1_CRR_DiadFit"
$ws.Range("B8").Value = "K21-64-FI#2_CRR_DiadFit"
$ws.Range("B9").Value = "K21-64-FI#8_CRR_DiadFit"

# Step 3: Update numeric result columns with refreshed CRR DiadFit values
# Row 2
$ws.Range("C2").Value = 103.0582439390753
$ws.Range("D2").Value = 1286.035116222911
$ws.Range("E2").Value = 592.2656617058601
$ws.Range("F2").Value = 1286.035116222911
$ws.Range("G2").Value = 607.8416133065451
$ws.Range("H2").Value = 0.433740954235018
$ws.Range("I2").Value = 3.692074287964963
$ws.Range("J2").Value = 0.3098621332226934
$ws.Range("K2").Value = 0.8674819084700359
$ws.Range("M2").Value = 1389.093360161986
$ws.Range("N2").Value = 968.2088497472205
$ws.Range("O2").Value = 1389.093360161986
$ws.Range("P2").Value = 904.3471424407591
$ws.Range("Q2").Value = 0.3954978032866402
$ws.Range("S2").Value = 3.743551243207165
$ws.Range("T2").Value = 0.3078957203892218
$ws.Range("U2").Value = 0.7909956065732805
# Row 3
$ws.Range("C3").Value = 103.0106887260733
$ws.Range("D3").Value = 1286.157462692443
$ws.Range("E3").Value = 227.6699200823184
$ws.Range("F3").Value = 1286.157462692443
$ws.Range("G3").Value = 201.5856524349188
$ws.Range("H3").Value = 0.4094465788269486
$ws.Range("I3").Value = 3.372069087053057
$ws.Range("J3").Value = 0.04144716748704053
$ws.Range("K3").Value = 0.8188931576538973
$ws.Range("M3").Value = 1389.168151418517
$ws.Range("N3").Value = 374.1384632773925
$ws.Range("O3").Value = 1389.168151418517
$ws.Range("P3").Value = 344.0597849544872
$ws.Range("Q3").Value = 0.3881542579594598
$ws.Range("S3").Value = 3.622729561803206
$ws.Range("T3").Value = 0.3189634016498229
$ws.Range("U3").Value = 0.7763085159189197
# Row 4
$ws.Range("C4").Value = 103.0463328358885
$ws.Range("D4").Value = 1286.074091074354
$ws.Range("E4").Value = 570.4088056961807
$ws.Range("F4").Value = 1286.074091074354
$ws.Range("G4").Value = 582.2966213040652
$ws.Range("H4").Value = 0.4222122028257283
$ws.Range("I4").Value = 3.767917484139669
$ws.Range("J4").Value = 0.3700715172096811
$ws.Range("K4").Value = 0.8444244056514566
$ws.Range("M4").Value = 1389.120423910243
$ws.Range("N4").Value = 888.6857868618297
$ws.Range("O4").Value = 1389.120423910243
$ws.Range("P4").Value = 839.6197158549975
$ws.Range("Q4").Value = 0.3996293612689302
$ws.Range("S4").Value = 3.922465505516284
$ws.Range("T4").Value = 0.3102593276707911
$ws.Range("U4").Value = 0.7992587225378603
# Row 5
$ws.Range("C5").Value = 103.3102715839677
$ws.Range("D5").Value = 1285.431835131392
$ws.Range("E5").Value = 65.96206751084628
$ws.Range("F5").Value = 1285.431835131392
$ws.Range("G5").Value = 81.92720526356536
$ws.Range("H5").Value = 0.5776929316504018
$ws.Range("I5").Value = 3.145730471303722
$ws.Range("J5").Value = 0.0000000001111621350524672
$ws.Range("K5").Value = 1.155385863300804
$ws.Range("M5").Value = 1388.74210671536
$ws.Range("N5").Value = 129.447908116705
$ws.Range("O5").Value = 1388.74210671536
$ws.Range("P5").Value = 167.5655996578714
$ws.Range("Q5").Value = 0.4790586387716631
$ws.Range("S5").Value = 2.390976642247799
$ws.Range("T5").Value = 0.6659269192991838
$ws.Range("U5").Value = 0.9581172775433262
# Row 6
$ws.Range("C6").Value = 103.3046133559294
$ws.Range("D6").Value = 1285.435763681008
$ws.Range("E6").Value = 607.8097089901272
$ws.Range("F6").Value = 1285.435763681008
$ws.Range("G6").Value = 895.8526333973542
$ws.Range("H6").Value = 0.5682520214930085
$ws.Range("I6").Value = 3.89174552535069
$ws.Range("J6").Value = 0.5562889388715584
$ws.Range("K6").Value = 1.136504042986017
$ws.Range("M6").Value = 1388.740377036938
$ws.Range("N6").Value = 1147.746984390957
$ws.Range("O6").Value = 1388.740377036938
$ws.Range("P6").Value = 1426.205823907998
$ws.Range("Q6").Value = 0.4810902539781914
$ws.Range("S6").Value = 4.056062987517096
$ws.Range("T6").Value = 0.5486866371465465
$ws.Range("U6").Value = 0.9621805079563829
# Row 7
$ws.Range("C7").Value = 103.3410505842598
$ws.Range("D7").Value = 1285.310757759881
$ws.Range("E7").Value = 125.5499373025313
$ws.Range("F7").Value = 1285.310757759881
$ws.Range("G7").Value = 210.7357224435302
$ws.Range("H7").Value = 0.6314607169829891
$ws.Range("I7").Value = 1.840530651270573
$ws.Range("J7").Value = 0.6219951783469195
$ws.Range("K7").Value = 1.262921433965978
$ws.Range("M7").Value = 1388.651808344141
$ws.Range("N7").Value = 235.7510865516424
$ws.Range("O7").Value = 1388.651808344141
$ws.Range("P7").Value = 317.8537635808829
$ws.Range("Q7").Value = 0.5313379959799696
$ws.Range("S7").Value = 1.767705392637869
$ws.Range("T7").Value = 0.5055209721595666
$ws.Range("U7").Value = 1.062675991959939
# Row 8
$ws.Range("C8").Value = 103.3010012985155
$ws.Range("D8").Value = 1285.40763681876
$ws.Range("E8").Value = 519.6086468321023
$ws.Range("F8").Value = 1285.40763681876
$ws.Range("G8").Value = 812.1829116361897
$ws.Range("H8").Value = 0.563375644419132
$ws.Range("I8").Value = 3.854109858933217
$ws.Range("J8").Value = 0.7243924632894531
$ws.Range("K8").Value = 1.126751288838264
$ws.Range("M8").Value = 1388.708638117275
$ws.Range("N8").Value = 938.3906529393677
$ws.Range("O8").Value = 1388.708638117275
$ws.Range("P8").Value = 1219.109980894617
$ws.Range("Q8").Value = 0.5075643996849418
$ws.Range("S8").Value = 3.674297197062103
$ws.Range("T8").Value = 0.5250338275249247
$ws.Range("U8").Value = 1.015128799369884
# Row 9
$ws.Range("C9").Value = 103.2888789760739
$ws.Range("D9").Value = 1285.392843382338
$ws.Range("E9").Value = 158.6481982128894
$ws.Range("F9").Value = 1285.392843382338
$ws.Range("G9").Value = 249.5983765892263
$ws.Range("H9").Value = 0.6197343560520949
$ws.Range("I9").Value = 1.485895158892154
$ws.Range("J9").Value = 0.5025210576443586
$ws.Range("K9").Value = 1.23946871210419
$ws.Range("M9").Value = 1388.681722358412
$ws.Range("N9").Value = 296.130691693525
$ws.Range("O9").Value = 1388.681722358412
$ws.Range("P9").Value = 390.8443694985737
$ws.Range("Q9").Value = 0.5144726707713018
$ws.Range("S9").Value = 1.934053969287134
$ws.Range("T9").Value = 0.5316083371779272
$ws.Range("U9").Value = 1.028945341542604

Write-Host "Done applying CRR DiadFit updates"
